$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 356, pushing the existing rows 356-406 down to 357-407.
$ws.Rows("356:356").Insert()

# Populate the newly inserted row 356 with the new weekly price record.
$ws.Cells.Item(356, 1).Value = 10
$ws.Cells.Item(356, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(356, 3).Value = "La Araucanía"
$ws.Cells.Item(356, 4).Value = 45131
$ws.Cells.Item(356, 5).Value = 9
$ws.Cells.Item(356, 6).Value = 100112039
$ws.Cells.Item(356, 7).Value = "Ciboulette"
$ws.Cells.Item(356, 8).Value = "Sin especificar"
$ws.Cells.Item(356, 9).Value = "Primera"
$ws.Cells.Item(356, 10).Value = 65
$ws.Cells.Item(356, 11).Value = 7000
$ws.Cells.Item(356, 12).Value = 7000
$ws.Cells.Item(356, 13).Value = 7000
$ws.Cells.Item(356, 14).Value = "$/docena de atados"
$ws.Cells.Item(356, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(356, 16).Value = 2333
$ws.Cells.Item(356, 17).Value = 3
$ws.Cells.Item(356, 18).Value = "Hortaliza"
